# Update cryptos list (Price / Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.180.71"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "2.629.33"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'518.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "

$ws.Range("D6").Value = "'148.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.10%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -3.78%  "

$ws.Range("D9").Value = "2.635.43"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("E10").Value = "  -5.31%  "

$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").Value = "3.088.01"
$ws.Range("E14").Value = "  +1.14%  "

$ws.Range("D15").Value = "60.161.46"
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("E16").Value = "  -2.29%  "

$ws.Range("E17").Value = "  -1.89%  "

$ws.Range("D18").Value = "2.620.39"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "'4.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.30%  "

$ws.Range("D20").Value = "'341.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.05%  "

$ws.Range("D21").Value = "'10.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").Value = "'61.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  -4.10%  "

$ws.Range("D28").Value = "0.0₃0809"
$ws.Range("E28").Value = "  -4.75%  "

$ws.Range("D29").Value = "'7.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.34%  "

$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").Value = "'5.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.62%  "

$ws.Range("E33").Value = "  -2.13%  "

$ws.Range("D34").Value = "'150.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("D35").Value = "'3.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.41%  "

$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("E37").Value = "  -5.52%  "

$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "'36.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("E40").Value = "  -4.82%  "

$ws.Range("E41").Value = "  -4.27%  "

$ws.Range("D42").Value = "'290.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.16%  "

$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("D44").Value = "'0.0999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "

$ws.Range("E46").Value = "  -2.55%  "

$ws.Range("D47").Value = "'19.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "

$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("D50").Value = "'4.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.24%  "

$ws.Range("D51").Value = "1.957.58"
$ws.Range("E51").Value = "  +0.06%  "

